# Apply the parameter-file update described in the commit
# "changes to param file and solver -100 to 100":
#   - Add two new zero-valued parameter columns (I4, J4) to the bottom
#     parameter row on Sheet1.
#   - Leave the active selection on I6, matching the saved cursor
#     position recorded in the workbook after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0

$ws.Range("I6").Select()
